$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "On profession" blurb: replace the description text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    ": Passionated with coding, performing experiments, learning and combining new technologies.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": Programming, researching for my thesis, learning about new technologies, such as Blockchain and Quantum Encryption.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "On spare-time" -> "In spare-time", and update its description text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "On spare-time",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In spare-time",
    2) | Out-Null

$d.Content.Find.Execute(
    ": Enjoying sports, reading books, and going crazy for travelling and food.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": Doing sports, reading books, and for travelling .",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) In the "Web-based device reservation system for JyvSecTec" paragraph,
#    insert a space between the ";" and the "Link" hyperlink that follows it.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Web-based device reservation system for JyvSecTec*") {
        $semiIdx = $t.IndexOf(";")
        $insertPos = $p.Range.Start + $semiIdx + 1
        $ins = $d.Range($insertPos, $insertPos)
        $ins.InsertAfter(" ")
        break
    }
}

# ---------------------------------------------------------------------------
# 4) In the "Thesis title: Software Defined Networking" paragraph, drop the
#    trailing " Link" hyperlink so the paragraph simply ends with ";".
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Thesis title: Software Defined Networking*") {
        $semiIdx = $t.IndexOf(";")
        $delStart = $p.Range.Start + $semiIdx + 1
        $delEnd = $p.Range.End - 1
        if ($delEnd -gt $delStart) {
            $delRange = $d.Range($delStart, $delEnd)
            $delRange.Delete()
        }
        break
    }
}
